$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 194-195, pushing the existing data (rows 194-292)
# down to rows 196-294.
$ws.Range("A194:R195").EntireRow.Insert()

# New row 194 - "Primera" quality observation dated 2021-11-18 (serial 44518)
$ws.Range("A194").Value = 9
$ws.Range("B194").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C194").Value = "Metropolitana"
$ws.Range("D194").Value = 44518
$ws.Range("D194").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E194").Value = 13
$ws.Range("F194").Value = 100112012
$ws.Range("G194").Value = "Espinaca"
$ws.Range("H194").Value = "Sin especificar"
$ws.Range("I194").Value = "Primera"
$ws.Range("J194").Value = 250
$ws.Range("K194").Value = 6000
$ws.Range("L194").Value = 7000
$ws.Range("M194").Value = 6500
$ws.Range("N194").Value = "$/cuna 10 kilos"
$ws.Range("O194").Value = "Provincia de Chacabuco"
$ws.Range("P194").Value = 650
$ws.Range("Q194").Value = 10
$ws.Range("R194").Value = "Hortaliza"

# New row 195 - "Segunda" quality observation dated 2021-11-18 (serial 44518)
$ws.Range("A195").Value = 9
$ws.Range("B195").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C195").Value = "Metropolitana"
$ws.Range("D195").Value = 44518
$ws.Range("D195").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E195").Value = 13
$ws.Range("F195").Value = 100112012
$ws.Range("G195").Value = "Espinaca"
$ws.Range("H195").Value = "Sin especificar"
$ws.Range("I195").Value = "Segunda"
$ws.Range("J195").Value = 106
$ws.Range("K195").Value = 5000
$ws.Range("L195").Value = 5000
$ws.Range("M195").Value = 5000
$ws.Range("N195").Value = "$/cuna 10 kilos"
$ws.Range("O195").Value = "Provincia de Chacabuco"
$ws.Range("P195").Value = 500
$ws.Range("Q195").Value = 10
$ws.Range("R195").Value = "Hortaliza"
